$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100; existing rows 100-142 shift down to 101-143.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new weekly price record.
$ws.Cells.Item(100, 1).Value = 3
$ws.Cells.Item(100, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(100, 3).Value = "Coquimbo"
$ws.Cells.Item(100, 4).Value = 44460
$ws.Cells.Item(100, 5).Value = 5
$ws.Cells.Item(100, 6).Value = 100112001
$ws.Cells.Item(100, 7).Value = "Berenjena"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 100
$ws.Cells.Item(100, 11).Value = 9500
$ws.Cells.Item(100, 12).Value = 10000
$ws.Cells.Item(100, 13).Value = 9700
$ws.Cells.Item(100, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(100, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(100, 16).Value = 162
$ws.Cells.Item(100, 17).Value = 60
$ws.Cells.Item(100, 18).Value = "Hortaliza"
